$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.414.24"
$ws.Range("E2").Value = "  -3.94%  "
$ws.Range("D3").Value = "2.649.14"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'522.24"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").Value = "'145.09"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'0.573"
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").Value = "'6.77"
$ws.Range("E9").Value = "  +3.91%  "
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").Value = "3.115.86"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").Value = "58.405.40"
$ws.Range("E14").Value = "  -3.88%  "
$ws.Range("D15").Value = "'20.97"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").Value = "'0.0000137"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "2.645.87"
$ws.Range("E17").Value = "  -15.03%  "
$ws.Range("D18").Value = "'339.03"
$ws.Range("D19").Value = "'4.40"
$ws.Range("E19").Value = "  -2.93%  "
$ws.Range("D20").Value = "'10.46"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").Value = "'6.33"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "'64.44"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "'0.167"
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").Value = "0.0₃0800"
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").Value = "'7.15"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").Value = "'6.70"
$ws.Range("E29").Value = "  -3.02%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").Value = "'152.47"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").Value = "'18.89"
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  -5.20%  "
$ws.Range("D36").Value = "'0.911"
$ws.Range("E36").Value = "  -4.48%  "
$ws.Range("D37").Value = "'0.871"
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").Value = "'36.80"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "'1.45"
$ws.Range("E39").Value = "  -4.88%  "
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'0.609"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "'274.49"
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("D44").Value = "'0.0972"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("D45").Value = "'19.46"
$ws.Range("D46").Value = "'0.0537"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("D48").Value = "2.042.99"
$ws.Range("E48").Value = "  -4.81%  "
$ws.Range("D49").Value = "'4.70"
$ws.Range("E49").Value = "  -4.39%  "
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("D51").Value = "'18.34"
$ws.Range("E51").Value = "  -3.74%  "
